# engine_node: fix color + add engine_node route for xl_json_converter
# Collapse the row-1/row-3 "list expansion" (columns E:G) back down to a
# single JSON-ish string value in D3, matching the non-expanded layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the now-unused trailing cells (both content and formatting) so the
# sheet's used range shrinks from A1:G3 down to A1:D3.
$ws.Range("E1:G1").Clear()
$ws.Range("E3:G3").Clear()

# D3 previously held the first element (0) of the expanded list; it now
# holds the whole list rendered as a single string.
$ws.Range("D3").Value = "[0, 1, 2, 3]"
